# Insert a new "fuel" worksheet right after "asymmetric_sgen" (i.e. before "ext_grid"),
# populate it with gen_type/index/fuel data, and make it the active sheet/tab —
# matching the target OOXML diff.

$wb = $excel.ActiveWorkbook

$anchor = $wb.Worksheets.Item("asymmetric_sgen")
$fuel = $wb.Worksheets.Add([System.Type]::Missing, $anchor)
$fuel.Name = "fuel"

# Header row (bold, size 12)
$fuel.Range("B1").Value = "gen_type"
$fuel.Range("C1").Value = "index"
$fuel.Range("D1").Value = "fuel"
$fuel.Range("B1:D1").Font.Bold = $true
$fuel.Range("B1:D1").Font.Size = 12

# Data rows: 8 sgens, all fuelled by solar
for ($i = 0; $i -lt 8; $i++) {
    $r = $i + 2
    $fuel.Cells.Item($r, 1).Value = $i
    $fuel.Cells.Item($r, 2).Value = "sgen"
    $fuel.Cells.Item($r, 3).Value = $i
    $fuel.Cells.Item($r, 4).Value = "solar"
}

# Leave the cursor where the author left it and make "fuel" the active/visible tab.
[void]$fuel.Range("M14").Select()
$fuel.Activate()
